$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '91.832.72'
$ws.Cells.Item(2, 5).Value = '  +0.74%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.099.31'
$ws.Cells.Item(3, 5).Value = '  -1.81%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.10%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '233.31'
$ws.Cells.Item(5, 5).Value = '  -2.71%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '611.69'
$ws.Cells.Item(6, 5).Value = '  -1.18%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.10'
$ws.Cells.Item(7, 5).Value = '  -2.34%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.385'
$ws.Cells.Item(8, 5).Value = '  +2.79%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '3.094.23'
$ws.Cells.Item(10, 5).Value = '  -1.95%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.763'
$ws.Cells.Item(11, 5).Value = '  +3.57%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.198'
$ws.Cells.Item(12, 5).Value = '  -2.86%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000243'
$ws.Cells.Item(13, 5).Value = '  -1.96%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '91.993.96'
$ws.Cells.Item(14, 5).Value = '  +1.29%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Toncoin'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.40'
$ws.Cells.Item(15, 5).Value = '  -4.26%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'Avalanche'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '33.56'
$ws.Cells.Item(16, 5).Value = '  -4.16%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.673.36'
$ws.Cells.Item(17, 5).Value = '  -1.88%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.092.00'
$ws.Cells.Item(18, 5).Value = '  -2.50%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.78'
$ws.Cells.Item(19, 5).Value = '  +1.99%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.44'
$ws.Cells.Item(20, 5).Value = '  -3.85%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.83'
$ws.Cells.Item(21, 5).Value = '  -2.07%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '437.81'
$ws.Cells.Item(22, 5).Value = '  -4.32%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.11'
$ws.Cells.Item(23, 5).Value = '  -0.47%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.0000194'
$ws.Cells.Item(24, 5).Value = '  -5.84%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '5.63'
$ws.Cells.Item(25, 5).Value = '  -6.25%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '85.64'
$ws.Cells.Item(26, 5).Value = '  -3.74%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.43'
$ws.Cells.Item(27, 5).Value = '  -4.42%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'WrappedeETH'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(28, 4).Value = '3.271.99'
$ws.Cells.Item(28, 5).Value = '  -1.85%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Dai'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.998'
$ws.Cells.Item(29, 5).Value = '  -0.16%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Hedera'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.131'
$ws.Cells.Item(30, 5).Value = '  -15.25%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Cronos'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.176'
$ws.Cells.Item(31, 5).Value = '  +4.88%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.236'
$ws.Cells.Item(32, 5).Value = '  +0.28%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.00'
$ws.Cells.Item(33, 5).Value = '  +0.24%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '9.02'
$ws.Cells.Item(34, 5).Value = '  -5.06%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'RenderToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '7.78'
$ws.Cells.Item(35, 5).Value = '  +2.15%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Kaspa'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.158'
$ws.Cells.Item(36, 5).Value = '  -8.85%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'EthereumClassic'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '25.57'
$ws.Cells.Item(37, 5).Value = '  -3.77%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.90'
$ws.Cells.Item(38, 5).Value = '  +1.31%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'PancakeSwap'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.88'
$ws.Cells.Item(39, 5).Value = '  -4.01%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '23.80'
$ws.Cells.Item(40, 5).Value = '  +7.57%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.438'
$ws.Cells.Item(41, 5).Value = '  -1.79%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Fetch.AI'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.27'
$ws.Cells.Item(42, 5).Value = '  -5.01%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Bittensor'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '466.49'
$ws.Cells.Item(43, 5).Value = '  -7.01%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'dogwifhat'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.23'
$ws.Cells.Item(44, 5).Value = '  -8.06%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'USDe'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  -0.04%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Monero'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '162.94'
$ws.Cells.Item(46, 5).Value = '  +4.66%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'ARBITRUM'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.680'
$ws.Cells.Item(47, 5).Value = '  -4.85%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Stacks'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.84'
$ws.Cells.Item(48, 5).Value = '  -4.13%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'ImmutableX'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.34'
$ws.Cells.Item(49, 5).Value = '  -2.84%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0325'
$ws.Cells.Item(50, 5).Value = '  -0.02%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'OKB'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '43.85'
$ws.Cells.Item(51, 5).Value = '  -0.63%  '
